# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7146
$ws1.Range("F6").Value = 558
$ws1.Range("F12").Value = 212
$ws1.Range("F14").Value = 456
$ws1.Range("F19").Value = 3715
$ws1.Range("F21").Value = 247
$ws1.Range("F25").Value = 2365
$ws1.Range("F27").Value = 284
$ws1.Range("F29").Value = 2
$ws1.Range("F30").Value = 39
$ws1.Range("F36").Value = 1387
$ws1.Range("F37").Value = 128

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7146
$ws4.Range("F7").Value = 558
$ws4.Range("F13").Value = 212
$ws4.Range("F15").Value = 456
$ws4.Range("F20").Value = 3715
$ws4.Range("F22").Value = 247
$ws4.Range("F26").Value = 2365
$ws4.Range("F28").Value = 284
$ws4.Range("F30").Value = 2
$ws4.Range("F31").Value = 39
$ws4.Range("F37").Value = 1387
$ws4.Range("F38").Value = 128

$wb.Save()
